# issue #5: add legislator_id, name, date into dataframe
#
# Adds three new trailing columns (date, legislator_name, legislator_id) to
# the "股票" (stocks) worksheet, filling every existing data row with the
# filing date, legislator name and legislator id for this report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

$legislatorName = "蔡正元"
$legislatorId = 966
$filingDate = "2012-04-18"

# --- Header row (row 1): copy the look of the existing header cell (G1) and
#     overwrite its text with the new column names. ---
$headerSrc = $ws.Cells.Item(1, 7)

$h1 = $ws.Cells.Item(1, 8)
$headerSrc.Copy($h1)
$h1.Value = "date"

$i1 = $ws.Cells.Item(1, 9)
$headerSrc.Copy($i1)
$i1.Value = "legislator_name"

$j1 = $ws.Cells.Item(1, 10)
$headerSrc.Copy($j1)
$j1.Value = "legislator_id"

# --- Data rows (rows 2-13): fill H/I/J. A scratch cell well outside the used
#     range is used to coerce the date string into the workbook as literal
#     text (avoiding automatic date-serial conversion) before its value is
#     pasted into place; cell formatting is copied separately from a normal
#     data cell on the same row so the new cells match their neighbours. ---
$scratch = $ws.Cells.Item(500, 500)

$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
for ($r = 2; $r -le $lastRow; $r++) {
    $fmtSrc = $ws.Cells.Item($r, 3)

    $dateCell = $ws.Cells.Item($r, 8)
    $scratch.Formula = "=""" + $filingDate + """"
    $scratch.Copy()
    $dateCell.PasteSpecial(-4163)
    $fmtSrc.Copy()
    $dateCell.PasteSpecial(-4122)

    $nameCell = $ws.Cells.Item($r, 9)
    $fmtSrc.Copy($nameCell)
    $nameCell.Value = $legislatorName

    $idCell = $ws.Cells.Item($r, 10)
    $fmtSrc.Copy($idCell)
    $idCell.Value = $legislatorId
}

$scratch.ClearContents()
